$wb = $excel.ActiveWorkbook

$wsALC = $wb.Worksheets.Item("ALC")
$wsARM = $wb.Worksheets.Item("ARM")
$wsBSM = $wb.Worksheets.Item("BSM")
$wsCRP = $wb.Worksheets.Item("CRP")
$wsCUL = $wb.Worksheets.Item("CUL")
$wsGSM = $wb.Worksheets.Item("GSM")
$wsLTW = $wb.Worksheets.Item("LTW")
$wsWVR = $wb.Worksheets.Item("WVR")

# ALC sheet updates
$wsALC.Range("H86").Value = 11143648
$wsALC.Range("I86").Value = 5380.5713
$wsALC.Range("K86").Value = 5380.5713
$wsALC.Range("M86").Value = -4257.5713
$wsALC.Range("H89").Value = 11143648
$wsALC.Range("I89").Value = 5380.5713
$wsALC.Range("K89").Value = 26902.8565
$wsALC.Range("M89").Value = -21286.8565
$wsALC.Range("H124").Value = 54999
$wsALC.Range("J124").Value = 54999
$wsALC.Range("L124").Value = 54999
$wsALC.Range("H130").Value = 61332.668
$wsALC.Range("J130").Value = 61332.668
$wsALC.Range("L130").Value = 61332.668
$wsALC.Range("N130").Value = -71372.66800000001
$wsALC.Range("H132").Value = 3898.95
$wsALC.Range("I132").Value = 4607.032
$wsALC.Range("J132").Value = 3142.0344
$wsALC.Range("K132").Value = 13821.096
$wsALC.Range("L132").Value = 9426.1032
$wsALC.Range("M132").Value = -11291.096
$wsALC.Range("N132").Value = -14486.1032
$wsALC.Range("H135").Value = 38462690
$wsALC.Range("I135").Value = 43479460
$wsALC.Range("J135").Value = 773.6667
$wsALC.Range("K135").Value = 391315140
$wsALC.Range("L135").Value = 6963.0003
$wsALC.Range("M135").Value = -391312605
$wsALC.Range("N135").Value = -12033.0003
$wsALC.Range("H137").Value = 71434010
$wsALC.Range("I137").Value = 200002800
$wsALC.Range("J137").Value = 6899.8887
$wsALC.Range("K137").Value = 600008400
$wsALC.Range("L137").Value = 20699.6661
$wsALC.Range("M137").Value = -600005850
$wsALC.Range("N137").Value = -25799.6661

# ARM sheet updates
$wsARM.Range("H61").Value = 2724.6365
$wsARM.Range("I61").Value = 1809
$wsARM.Range("K61").Value = 1809
$wsARM.Range("M61").Value = -1597
$wsARM.Range("H74").Value = 2781.35
$wsARM.Range("I74").Value = 2806.7222
$wsARM.Range("K74").Value = 2806.7222
$wsARM.Range("M74").Value = -1932.7222
$wsARM.Range("H77").Value = 2781.35
$wsARM.Range("I77").Value = 2806.7222
$wsARM.Range("K77").Value = 14033.611
$wsARM.Range("M77").Value = -9665.611000000001
$wsARM.Range("H97").Value = 291.72726
$wsARM.Range("I97").Value = 300.9
$wsARM.Range("J97").Value = 200
$wsARM.Range("K97").Value = 300.9
$wsARM.Range("L97").Value = 200
$wsARM.Range("M97").Value = 195.1
$wsARM.Range("N97").Value = -1192
$wsARM.Range("H136").Value = 2724.6365
$wsARM.Range("I136").Value = 1809
$wsARM.Range("K136").Value = 5427
$wsARM.Range("M136").Value = -2877

# BSM sheet updates
$wsBSM.Range("H20").Value = 4294.9287
$wsBSM.Range("I20").Value = 3856.077
$wsBSM.Range("K20").Value = 3856.077
$wsBSM.Range("M20").Value = -3609.077
$wsBSM.Range("H86").Value = 33335538
$wsBSM.Range("I86").Value = 45456428
$wsBSM.Range("K86").Value = 45456428
$wsBSM.Range("M86").Value = -45455305
$wsBSM.Range("H89").Value = 33335538
$wsBSM.Range("I89").Value = 45456428
$wsBSM.Range("K89").Value = 227282140
$wsBSM.Range("M89").Value = -227276524
$wsBSM.Range("H105").Value = 2776.6365
$wsBSM.Range("I105").Value = 2597.4
$wsBSM.Range("K105").Value = 2597.4
$wsBSM.Range("M105").Value = -850.4000000000001
$wsBSM.Range("H130").Value = 60000
$wsBSM.Range("J130").Value = 60000
$wsBSM.Range("L130").Value = 60000
$wsBSM.Range("N130").Value = -70040

# CRP sheet updates
$wsCRP.Range("H62").Value = 6742
$wsCRP.Range("I62").Value = 5657.4707
$wsCRP.Range("J62").Value = 9814.833000000001
$wsCRP.Range("K62").Value = 5657.4707
$wsCRP.Range("L62").Value = 9814.833000000001
$wsCRP.Range("M62").Value = -5033.4707
$wsCRP.Range("N62").Value = -11062.833
$wsCRP.Range("H65").Value = 6742
$wsCRP.Range("I65").Value = 5657.4707
$wsCRP.Range("J65").Value = 9814.833000000001
$wsCRP.Range("K65").Value = 28287.3535
$wsCRP.Range("L65").Value = 49074.165
$wsCRP.Range("M65").Value = -25167.3535
$wsCRP.Range("N65").Value = -55314.165
$wsCRP.Range("H68").Value = 52499.168
$wsCRP.Range("J68").Value = 52499.168
$wsCRP.Range("L68").Value = 52499.168
$wsCRP.Range("N68").Value = -53997.168
$wsCRP.Range("H71").Value = 52499.168
$wsCRP.Range("J71").Value = 52499.168
$wsCRP.Range("L71").Value = 157497.504
$wsCRP.Range("N71").Value = -164985.504
$wsCRP.Range("H107").Value = 776.0625
$wsCRP.Range("I107").Value = 835.6667
$wsCRP.Range("J107").Value = 597.25
$wsCRP.Range("K107").Value = 835.6667
$wsCRP.Range("L107").Value = 597.25
$wsCRP.Range("M107").Value = 1084.3333
$wsCRP.Range("N107").Value = -4437.25
$wsCRP.Range("H130").Value = 51073.332
$wsCRP.Range("J130").Value = 51073.332
$wsCRP.Range("L130").Value = 51073.332
$wsCRP.Range("N130").Value = -61113.332
$wsCRP.Range("H132").Value = 3034.84
$wsCRP.Range("I132").Value = 2108.647
$wsCRP.Range("K132").Value = 6325.941
$wsCRP.Range("M132").Value = -3795.941

# CUL sheet updates
$wsCUL.Range("H107").Value = 1419.8182
$wsCUL.Range("J107").Value = 0
$wsCUL.Range("L107").Value = 0

# GSM sheet updates
$wsGSM.Range("H51").Value = 29999
$wsGSM.Range("J51").Value = 29999
$wsGSM.Range("L51").Value = 29999
$wsGSM.Range("N51").Value = -31017
$wsGSM.Range("H80").Value = 3629.7646
$wsGSM.Range("J80").Value = 3728.6
$wsGSM.Range("L80").Value = 3728.6
$wsGSM.Range("N80").Value = -5724.6
$wsGSM.Range("H83").Value = 3629.7646
$wsGSM.Range("J83").Value = 3728.6
$wsGSM.Range("L83").Value = 18643
$wsGSM.Range("N83").Value = -28627
$wsGSM.Range("H102").Value = 3371.2
$wsGSM.Range("I102").Value = 2456.7273
$wsGSM.Range("K102").Value = 2456.7273
$wsGSM.Range("M102").Value = -834.7273
$wsGSM.Range("H124").Value = 54995.332
$wsGSM.Range("J124").Value = 54995.332
$wsGSM.Range("L124").Value = 54995.332
$wsGSM.Range("N124").Value = -64815.332
$wsGSM.Range("H128").Value = 49499.332
$wsGSM.Range("J128").Value = 49499.332
$wsGSM.Range("L128").Value = 49499.332
$wsGSM.Range("N128").Value = -59459.332

# LTW sheet updates
$wsLTW.Range("H16").Value = 12333.333
$wsLTW.Range("H68").Value = 8024.0557
$wsLTW.Range("I68").Value = 3448.75
$wsLTW.Range("J68").Value = 9331.286
$wsLTW.Range("K68").Value = 3448.75
$wsLTW.Range("L68").Value = 9331.286
$wsLTW.Range("M68").Value = -2699.75
$wsLTW.Range("N68").Value = -10829.286
$wsLTW.Range("H71").Value = 8024.0557
$wsLTW.Range("I71").Value = 3448.75
$wsLTW.Range("J71").Value = 9331.286
$wsLTW.Range("K71").Value = 17243.75
$wsLTW.Range("L71").Value = 46656.43
$wsLTW.Range("M71").Value = -13499.75
$wsLTW.Range("N71").Value = -54144.43
$wsLTW.Range("H132").Value = 5034.7744
$wsLTW.Range("I132").Value = 2640.2778
$wsLTW.Range("K132").Value = 7920.8334
$wsLTW.Range("M132").Value = -5390.8334

# WVR sheet updates
$wsWVR.Range("H43").Value = 18999
$wsWVR.Range("I43").Value = 18999
$wsWVR.Range("K43").Value = 18999
$wsWVR.Range("M43").Value = -18850

# Add new cell N124 on ALC
$wsALC.Range("N124").Value = -64819

# Remove cell N107 on CUL
$wsCUL.Range("N107").ClearContents()
